$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.152.94"

# Row 3
$ws.Range("D3").Value = "1.623.87"
$ws.Range("E3").Value = "  -1.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Formula = "'214.25"

# Row 6
$ws.Range("D6").Formula = "'0.523"
$ws.Range("E6").Value = "  +1.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -1.50%  "

# Row 9
$ws.Range("E9").Value = "  -0.17%  "

# Row 10
$ws.Range("E10").Value = "  +1.25%  "

# Row 11
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "1.613.78"
$ws.Range("E12").Value = "  -2.30%  "

# Row 13
$ws.Range("E13").Value = "  -0.43%  "

# Row 15
$ws.Range("D15").Value = "27.138.79"

# Row 16
$ws.Range("D16").Formula = "'64.57"
$ws.Range("E16").Value = "  -4.23%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0746"
$ws.Range("E17").Value = "  +0.47%  "

# Row 18
$ws.Range("D18").Formula = "'215.75"

# Row 19
$ws.Range("E19").Value = "  +0.01%  "

# Row 20
$ws.Range("D20").Formula = "'6.91"
$ws.Range("E20").Value = "  +0.36%  "

# Row 22
$ws.Range("D22").Formula = "'2.41"
$ws.Range("E22").Value = "  -6.45%  "

# Row 23
$ws.Range("E23").Value = "  -2.04%  "

# Row 24
$ws.Range("D24").Formula = "'148.13"
$ws.Range("E24").Value = "  +0.34%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("E26").Value = "  -3.20%  "

# Row 27
$ws.Range("E27").Value = "  -1.25%  "

# Row 28
$ws.Range("E28").Value = "  -1.14%  "

# Row 29
$ws.Range("E29").Value = "  -0.64%  "

# Row 31
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("E32").Value = "  -1.06%  "

# Row 33
$ws.Range("D33").Value = "1.343.47"
$ws.Range("E33").Value = "  +5.24%  "

# Row 34
$ws.Range("E34").Value = "  -0.54%  "

# Row 36
$ws.Range("D36").Formula = "'0.0178"

# Row 37
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
$ws.Range("D38").Formula = "'0.859"
$ws.Range("E38").Value = "  -0.24%  "

# Row 39
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("D40").Formula = "'0.802"

# Row 41
$ws.Range("D41").Formula = "'65.59"
$ws.Range("E41").Value = "  +5.84%  "

# Row 42
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("E43").Value = "  -1.30%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Formula = "'0.930"
$ws.Range("E44").Value = "  +38.75%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.760.67"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46
$ws.Range("D46").Formula = "'89.95"
$ws.Range("E46").Value = "  -2.19%  "

# Row 47
$ws.Range("E47").Value = "  +1.03%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Formula = "'0.100"
$ws.Range("E48").Value = "  +2.79%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Formula = "'0.0514"
$ws.Range("E49").Value = "  -0.39%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Formula = "'7.58"
$ws.Range("E50").Value = "  -0.95%  "

# Row 51
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Formula = "'1.00"
$ws.Range("E51").Value = "  -0.07%  "
